$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 808.4375
$ws.Range("I32").Value = 430.54544
$ws.Range("J32").Value = 1639.8
$ws.Range("K32").Value = 430.54544
$ws.Range("L32").Value = 1639.8
$ws.Range("M32").Value = -104.54544
$ws.Range("N32").Value = -2291.8
$ws.Range("H87").Value = 15391.963
$ws.Range("J87").Value = 15391.963
$ws.Range("L87").Value = 15391.963
$ws.Range("N87").Value = -17887.963
$ws.Range("H90").Value = 15391.963
$ws.Range("J90").Value = 15391.963
$ws.Range("L90").Value = 46175.889
$ws.Range("N90").Value = -58655.889
$ws.Range("H98").Value = 1368.7059
$ws.Range("I98").Value = 1391.6923
$ws.Range("J98").Value = 1294
$ws.Range("K98").Value = 1391.6923
$ws.Range("L98").Value = 1294
$ws.Range("M98").Value = 106.3077000000001
$ws.Range("N98").Value = -4290
$ws.Range("H122").Value = 1368.7059
$ws.Range("I122").Value = 1391.6923
$ws.Range("J122").Value = 1294
$ws.Range("K122").Value = 4175.0769
$ws.Range("L122").Value = 3882
$ws.Range("M122").Value = -1725.0769
$ws.Range("N122").Value = -8782
$ws.Range("H129").Value = 904.8148
$ws.Range("J129").Value = 931.2
$ws.Range("L129").Value = 2793.6
$ws.Range("N129").Value = -12793.6
$ws.Range("H137").Value = 1255
$ws.Range("I137").Value = 1258.7368
$ws.Range("J137").Value = 1249.9286
$ws.Range("K137").Value = 3776.2104
$ws.Range("L137").Value = 3749.7858
$ws.Range("M137").Value = -1226.2104
$ws.Range("N137").Value = -8849.7858
$ws.Range("H138").Value = 4203.01
$ws.Range("J138").Value = 4900.364
$ws.Range("L138").Value = 14701.092
$ws.Range("N138").Value = -24981.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1499.4308
$ws.Range("I74").Value = 968.25
$ws.Range("J74").Value = 3624.1538
$ws.Range("K74").Value = 968.25
$ws.Range("L74").Value = 3624.1538
$ws.Range("M74").Value = -94.25
$ws.Range("N74").Value = -5372.1538
$ws.Range("H77").Value = 1499.4308
$ws.Range("I77").Value = 968.25
$ws.Range("J77").Value = 3624.1538
$ws.Range("K77").Value = 4841.25
$ws.Range("L77").Value = 18120.769
$ws.Range("M77").Value = -473.25
$ws.Range("N77").Value = -26856.769
$ws.Range("H97").Value = 649.63336
$ws.Range("I97").Value = 376.35294
$ws.Range("K97").Value = 376.35294
$ws.Range("M97").Value = 119.64706
$ws.Range("H122").Value = 2295.647
$ws.Range("I122").Value = 1417.6666
$ws.Range("K122").Value = 4252.9998
$ws.Range("M122").Value = -1802.9998
$ws.Range("H132").Value = 2585.575
$ws.Range("I132").Value = 2185.2693
$ws.Range("J132").Value = 3329
$ws.Range("K132").Value = 6555.8079
$ws.Range("L132").Value = 9987
$ws.Range("M132").Value = -4025.8079
$ws.Range("N132").Value = -15047

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 22711.2
$ws.Range("J52").Value = 22711.2
$ws.Range("L52").Value = 22711.2
$ws.Range("N52").Value = -23237.2
$ws.Range("H99").Value = 1947.8334
$ws.Range("I99").Value = 1184
$ws.Range("J99").Value = 2902.625
$ws.Range("K99").Value = 1184
$ws.Range("L99").Value = 2902.625
$ws.Range("M99").Value = 314
$ws.Range("N99").Value = -5898.625
$ws.Range("H121").Value = 22711.2
$ws.Range("J121").Value = 22711.2
$ws.Range("L121").Value = 22711.2
$ws.Range("N121").Value = -26205.2
$ws.Range("H132").Value = 45937.5
$ws.Range("J132").Value = 45937.5
$ws.Range("L132").Value = 45937.5
$ws.Range("N132").Value = -56057.5
$ws.Range("H134").Value = 1821.2916
$ws.Range("I134").Value = 1414.4242
$ws.Range("J134").Value = 2716.4
$ws.Range("K134").Value = 4243.2726
$ws.Range("L134").Value = 8149.200000000001
$ws.Range("M134").Value = -1708.2726
$ws.Range("N134").Value = -13219.2
$ws.Range("H140").Value = 58868
$ws.Range("J140").Value = 58868
$ws.Range("L140").Value = 58868
$ws.Range("N140").Value = -69228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11425.294
$ws.Range("I31").Value = 5743.143
$ws.Range("J31").Value = 15402.8
$ws.Range("K31").Value = 5743.143
$ws.Range("L31").Value = 15402.8
$ws.Range("M31").Value = -5448.143
$ws.Range("N31").Value = -15992.8
$ws.Range("H34").Value = 11425.294
$ws.Range("I34").Value = 5743.143
$ws.Range("J34").Value = 15402.8
$ws.Range("K34").Value = 5743.143
$ws.Range("L34").Value = 15402.8
$ws.Range("M34").Value = -5541.143
$ws.Range("N34").Value = -15806.8
$ws.Range("H62").Value = 5575.4165
$ws.Range("I62").Value = 6602.5
$ws.Range("J62").Value = 5370
$ws.Range("K62").Value = 6602.5
$ws.Range("L62").Value = 5370
$ws.Range("M62").Value = -5978.5
$ws.Range("N62").Value = -6618
$ws.Range("H65").Value = 5575.4165
$ws.Range("I65").Value = 6602.5
$ws.Range("J65").Value = 5370
$ws.Range("K65").Value = 33012.5
$ws.Range("L65").Value = 26850
$ws.Range("M65").Value = -29892.5
$ws.Range("N65").Value = -33090
$ws.Range("H107").Value = 491.5909
$ws.Range("I107").Value = 272.33334
$ws.Range("J107").Value = 1478.25
$ws.Range("K107").Value = 272.33334
$ws.Range("L107").Value = 1478.25
$ws.Range("M107").Value = 1647.66666
$ws.Range("N107").Value = -5318.25
$ws.Range("H138").Value = 40285.445
$ws.Range("J138").Value = 40285.445
$ws.Range("L138").Value = 40285.445
$ws.Range("N138").Value = -50565.445
$ws.Range("H140").Value = 74254
$ws.Range("J140").Value = 74254
$ws.Range("L140").Value = 74254
$ws.Range("N140").Value = -84614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 821.84
$ws.Range("J131").Value = 894.3855600000001
$ws.Range("L131").Value = 2683.15668
$ws.Range("N131").Value = -12763.15668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2197.6365
$ws.Range("I113").Value = 1810
$ws.Range("J113").Value = 2662.8
$ws.Range("K113").Value = 1810
$ws.Range("L113").Value = 2662.8
$ws.Range("M113").Value = 360
$ws.Range("N113").Value = -7002.8
$ws.Range("H126").Value = 2499.0527
$ws.Range("I126").Value = 2480
$ws.Range("J126").Value = 2512.9092
$ws.Range("K126").Value = 7440
$ws.Range("L126").Value = 7538.7276
$ws.Range("M126").Value = -4970
$ws.Range("N126").Value = -12478.7276
$ws.Range("H132").Value = 2629
$ws.Range("I132").Value = 2277.577
$ws.Range("J132").Value = 3644.2222
$ws.Range("K132").Value = 6832.731000000001
$ws.Range("L132").Value = 10932.6666
$ws.Range("M132").Value = -4302.731000000001
$ws.Range("N132").Value = -15992.6666
$ws.Range("H135").Value = 42000
$ws.Range("J135").Value = 42000
$ws.Range("L135").Value = 42000
$ws.Range("N135").Value = -52140
$ws.Range("H140").Value = 39600
$ws.Range("J140").Value = 39600
$ws.Range("L140").Value = 39600
$ws.Range("N140").Value = -49960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1269.1666
$ws.Range("I22").Value = 919.5
$ws.Range("J22").Value = 1706.25
$ws.Range("K22").Value = 919.5
$ws.Range("L22").Value = 1706.25
$ws.Range("M22").Value = -624.5
$ws.Range("N22").Value = -2296.25
$ws.Range("H27").Value = 1269.1666
$ws.Range("I27").Value = 919.5
$ws.Range("J27").Value = 1706.25
$ws.Range("K27").Value = 919.5
$ws.Range("L27").Value = 1706.25
$ws.Range("M27").Value = -812.5
$ws.Range("N27").Value = -1920.25
$ws.Range("H40").Value = 250750
$ws.Range("I40").Value = 250750
$ws.Range("K40").Value = 250750
$ws.Range("M40").Value = -250614
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1433.3334
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 1433.3334
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -892.3334
$ws.Range("N100").Value = -2682
$ws.Range("H127").Value = 70686
$ws.Range("J127").Value = 70686
$ws.Range("L127").Value = 70686
$ws.Range("N127").Value = -80606
$ws.Range("H133").Value = 46326
$ws.Range("J133").Value = 46326
$ws.Range("L133").Value = 46326
$ws.Range("N133").Value = -51386
$ws.Range("H136").Value = 12823852
$ws.Range("I136").Value = 3831.5881
$ws.Range("K136").Value = 11494.7643
$ws.Range("M136").Value = -8944.764299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 40000
$ws.Range("J42").Value = 40000
$ws.Range("L42").Value = 40000
$ws.Range("N42").Value = -40756
$ws.Range("H81").Value = 2825
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2825
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -40608
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360
